$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E to remain text (matches original inlineStr text pattern)
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "51.456.42"
$ws.Range("E2").Value = "  -1.25%  "

# Row 3
$ws.Range("D3").Value = "2.925.00"
$ws.Range("E3").Value = "  -0.48%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "350.54"
$ws.Range("E5").Value = "  -0.74%  "

# Row 6
$ws.Range("D6").Value = "107.08"
$ws.Range("E6").Value = "  -4.52%  "

# Row 7
$ws.Range("E7").Value = "  -1.41%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  -3.26%  "

# Row 10
$ws.Range("D10").Value = "37.75"
$ws.Range("E10").Value = "  -4.34%  "

# Row 11
$ws.Range("E11").Value = "  +1.28%  "

# Row 12
$ws.Range("D12").Value = "0.0847"
$ws.Range("E12").Value = "  -3.47%  "

# Row 13
$ws.Range("D13").Value = "18.90"
$ws.Range("E13").Value = "  -6.14%  "

# Row 14
$ws.Range("D14").Value = "3.377.75"
$ws.Range("E14").Value = "  -0.71%  "

# Row 15
$ws.Range("D15").Value = "7.52"
$ws.Range("E15").Value = "  -3.17%  "

# Row 16
$ws.Range("D16").Value = "2.921.12"
$ws.Range("E16").Value = "  -0.88%  "

# Row 17
$ws.Range("D17").Value = "0.961"
$ws.Range("E17").Value = "  -2.13%  "

# Row 18
$ws.Range("D18").Value = "51.396.12"
$ws.Range("E18").Value = "  -1.40%  "

# Row 19
$ws.Range("D19").Value = "3.41"
$ws.Range("E19").Value = "  +3.50%  "

# Row 20
$ws.Range("D20").Value = "7.39"
$ws.Range("E20").Value = "  -2.98%  "

# Row 21
$ws.Range("D21").Value = "13.40"
$ws.Range("E21").Value = "  -5.89%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0965"

# Row 23
$ws.Range("D23").Value = "68.86"
$ws.Range("E23").Value = "  -3.33%  "

# Row 24
$ws.Range("D24").Value = "260.25"
$ws.Range("E24").Value = "  -3.11%  "

# Row 25
$ws.Range("E25").Value = "  -3.14%  "

# Row 26
$ws.Range("E26").Value = "  -3.48%  "

# Row 27
$ws.Range("D27").Value = "26.38"
$ws.Range("E27").Value = "  -2.50%  "

# Row 28
$ws.Range("E28").Value = "  +0.11%  "

# Row 29
$ws.Range("D29").Value = "7.36"
$ws.Range("E29").Value = "  +1.52%  "

# Row 30
$ws.Range("E30").Value = "  +0.02%  "

# Row 31
$ws.Range("E31").Value = "  -3.83%  "

# Row 32
$ws.Range("D32").Value = "6.06"
$ws.Range("E32").Value = "  -0.80%  "

# Row 33
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "2.19"
$ws.Range("E33").Value = "  -3.05%  "

# Row 34
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "35.63"
$ws.Range("E34").Value = "  -3.91%  "

# Row 35
$ws.Range("D35").Value = "50.35"
$ws.Range("E35").Value = "  -5.14%  "

# Row 36
$ws.Range("D36").Value = "0.0427"
$ws.Range("E36").Value = "  -5.83%  "

# Row 37
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.06%  "

# Row 38
$ws.Range("D38").Value = "3.13"
$ws.Range("E38").Value = "  -7.50%  "

# Row 39
$ws.Range("D39").Value = "17.61"
$ws.Range("E39").Value = "  -5.79%  "

# Row 40
$ws.Range("E40").Value = "  -5.98%  "

# Row 41
$ws.Range("E41").Value = "  -1.91%  "

# Row 42
$ws.Range("E42").Value = "  -1.98%  "

# Row 43
$ws.Range("D43").Value = "22.21"
$ws.Range("E43").Value = "  -4.71%  "

# Row 44
$ws.Range("D44").Value = "119.62"
$ws.Range("E44").Value = "  +7.16%  "

# Row 45
$ws.Range("E45").Value = "  -3.29%  "

# Row 46
$ws.Range("D46").Value = "2.094.41"

# Row 47
$ws.Range("D47").Value = "3.31"
$ws.Range("E47").Value = "  -6.40%  "

# Row 48
$ws.Range("D48").Value = "2.30"
$ws.Range("E48").Value = "  -9.05%  "

# Row 49
$ws.Range("E49").Value = "  -4.01%  "

# Row 50
$ws.Range("D50").Value = "0.0333"
$ws.Range("E50").Value = "  -5.92%  "

# Row 51
$ws.Range("D51").Value = "0.906"
$ws.Range("E51").Value = "  -5.12%  "
